# Actualización automática 2025-11-24 08:30:09
#
# Applies the updated November sales figures across the three sheets:
#   - "VENTAS POR GRUPO"       : per-client / per-product-group sales
#   - "VENTA MENSUAL"          : per-client monthly sales (column F = noviembre)
#   - "CUMPLIMIENTO MENSUAL"   : per product-group budget compliance summary

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO  (sheet1) - new sales booked against product groups
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("D5").Value  = 146.53      # CARRION CARRION LESLY ANABE - 240X80 PORCELANATO
$wsGrupo.Range("D6").Value  = 8875.02     # CHASIQUIZA CAMPAÑA JOSE LUIS - 240X80 PORCELANATO
$wsGrupo.Range("I11").Value = 384.5       # JARAMILLO CARVAJAL NICOLAS ESTEBAN - LAVABOS
$wsGrupo.Range("M11").Value = 5809.48     # JARAMILLO CARVAJAL NICOLAS ESTEBAN - PORCELANATO
$wsGrupo.Range("M13").Value = 6343.23     # MEGAMAFERS S.A. - PORCELANATO
$wsGrupo.Range("K14").Value = 690.28      # MUÑOZ LOZA ROMMEL SEBASTIAN - PANELES DECORATIVOS

# Row 23 totals how many (of the 21) rows have a positive value per column.
$wsGrupo.Range("D23").Value = "3 de 21"
$wsGrupo.Range("I23").Value = "2 de 21"
$wsGrupo.Range("K23").Value = "1 de 21"

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL  (sheet2) - same clients, November ("noviembre") column
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F5").Value  = -1654.27
$wsMensual.Range("F6").Value  = 8875.02
$wsMensual.Range("F11").Value = 6193.98
$wsMensual.Range("F13").Value = 6343.23
$wsMensual.Range("F14").Value = 3015.09
$wsMensual.Range("F23").Value = 21283.68   # column total

# ---------------------------------------------------------------------
# 3) CUMPLIMIENTO MENSUAL  (sheet3) - budget compliance per product group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 9326.83
$wsCumpl.Range("E3").Value = -6762.83
$wsCumpl.Range("F3").Value = 3.637609204368175

# LAVABOS
$wsCumpl.Range("D7").Value = 428.6
$wsCumpl.Range("E7").Value = -45.20000000000005
$wsCumpl.Range("F7").Value = 1.117892540427752

# PANELES DECORATIVOS
$wsCumpl.Range("D10").Value = 690.28
$wsCumpl.Range("E10").Value = 697.72
$wsCumpl.Range("F10").Value = 0.4973198847262247

# PORCELANATO
$wsCumpl.Range("D12").Value = 10824.3
$wsCumpl.Range("E12").Value = 33593.7
$wsCumpl.Range("F12").Value = 0.2436917465892206

# TOTAL
$wsCumpl.Range("D14").Value = 21283.68
$wsCumpl.Range("E14").Value = 34115.79101170094
$wsCumpl.Range("F14").Value = 0.3841856178645581

# The longer negative value in E7 makes Excel widen column E (5th column)
# when the sheet is re-saved; match the resulting autofit width.
$wsCumpl.Range("E1").EntireColumn.ColumnWidth = 23.17
